# Add a "react_basic" style roster/score layout in columns G:H (names + scores)
# next to the existing API table, matching the project-ideas update described
# in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Names go into column G (as text, so they land in sharedStrings).
# The insertion order below reproduces the exact order the strings were
# originally typed in (rows 2-5 top-down, then row 10 back up through 6,
# then rows 11-26 top-down).
$ws.Range("G2").Value = "Harsha"
$ws.Range("G3").Value = "Manashvini"
$ws.Range("G4").Value = "Charan"
$ws.Range("G5").Value = "Yashaswini "
$ws.Range("G10").Value = "Nilima"
$ws.Range("G9").Value = "SriHari"
$ws.Range("G8").Value = "Vineeth"
$ws.Range("G7").Value = "Kushbu"
$ws.Range("G6").Value = "Sathvika"
$ws.Range("G11").Value = "Revanth"
$ws.Range("G12").Value = "HariPriya"
$ws.Range("G13").Value = "Sai Vamsi"
$ws.Range("G14").Value = "Supriya"
$ws.Range("G15").Value = "Rohith "
$ws.Range("G16").Value = "Shraviya"
$ws.Range("G17").Value = "Vivek"
$ws.Range("G18").Value = "Sanmuk"
$ws.Range("G19").Value = "Kumar"
$ws.Range("G20").Value = "Ramesh"
$ws.Range("G21").Value = "Thusar"
$ws.Range("G22").Value = "Pranitha"
$ws.Range("G23").Value = "Bhaskar "
$ws.Range("G24").Value = "Keerthana"
$ws.Range("G25").Value = "Poojitha"
$ws.Range("G26").Value = "Ramaraju"

# Matching numeric scores in column H.
$ws.Range("H2").Value = 80
$ws.Range("H3").Value = 80
$ws.Range("H4").Value = 60
$ws.Range("H5").Value = 60
$ws.Range("H6").Value = 60
$ws.Range("H7").Value = 20
$ws.Range("H8").Value = 30
$ws.Range("H9").Value = 40
$ws.Range("H10").Value = 50
$ws.Range("H11").Value = 95
$ws.Range("H12").Value = 80
$ws.Range("H13").Value = 50
$ws.Range("H14").Value = 60
$ws.Range("H15").Value = 20
$ws.Range("H16").Value = 20
$ws.Range("H17").Value = 80
$ws.Range("H18").Value = 95
$ws.Range("H19").Value = 40
$ws.Range("H20").Value = 70
$ws.Range("H21").Value = 50
$ws.Range("H22").Value = 30
$ws.Range("H23").Value = 30
$ws.Range("H24").Value = 40
$ws.Range("H25").Value = 80
$ws.Range("H26").Value = 60

# Widen column G to fit the new names (closest achievable width to the
# original 10.33203125 best-fit result).
$ws.Columns("G").ColumnWidth = 9.5

# Leave the selection on the last-entered cell, as in the saved workbook.
$ws.Range("G17").Select() | Out-Null
